$wb = $excel.ActiveWorkbook

# --- 1. Add a new "Player Info" sheet, placed before the existing "ODI Batting" sheet ---
$battingRef = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingRef)
$playerInfo.Name = "Player Info"

# NOTE: after Worksheets.Add(), previously obtained sheet references can point at
# the wrong sheet (the reference tracks position, and the new sheet took that
# position). Re-fetch the sheets we still need to edit, by name, now.
$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $playerInfo.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$idCell = $playerInfo.Cells.Item(2, 1)
$idCell.NumberFormat = "@"
$idCell.Value = "5662"
$idCell.ClearFormats()
$playerInfo.Cells.Item(2, 2).Value = "Romario Shepherd"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Medium Fast"

# --- 2. Update "ODI Batting" sheet: column D header + values (URL -> match code) ---
$battingRows = $batting.UsedRange.Rows.Count
$batting.Cells.Item(1, 4).Value = "MATCH_CODE"
for ($r = 2; $r -le $battingRows; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $link = $cell.Value2
    if ($link) {
        $code = $link -replace '^.*MatchCode=', ''
        $cell.NumberFormat = "@"
        $cell.Value = $code
        $cell.ClearFormats()
    }
}

# --- 3. Update "ODI Bowling" sheet: column B header + values (URL -> match code) ---
$bowlingRows = $bowling.UsedRange.Rows.Count
$bowling.Cells.Item(1, 2).Value = "MATCH_CODE"
for ($r = 2; $r -le $bowlingRows; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $link = $cell.Value2
    if ($link) {
        $code = $link -replace '^.*MatchCode=', ''
        $cell.NumberFormat = "@"
        $cell.Value = $code
        $cell.ClearFormats()
    }
}
